$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: remove the "recently updated" yellow-highlight style from an N-col
# date cell by pasting formats from a cell that already carries the plain
# (non-highlighted) date style, without disturbing the cell's own value.
# ---------------------------------------------------------------------------
function Clear-DateHighlight($cellRef) {
    $src = $ws.Range("N18")        # known to carry the plain (non-highlighted) date style
    $dst = $ws.Range($cellRef)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# --- Rows 10-12: JOLTS figures refreshed this cycle -> drop the highlight ---
Clear-DateHighlight("N10")
Clear-DateHighlight("N11")
Clear-DateHighlight("N12")

# --- Row 13: UI Initial Claims ---
$ws.Range("N13").Value = 46055
$ws.Range("Q13").Value = 227000
$ws.Range("R13").Value = 232000
$ws.Range("S13").Value = 209000
$ws.Range("T13").Value = 210000
$ws.Range("U13").Value = 199000

# --- Row 14: UI Continuing Claims ---
$ws.Range("N14").Value = 46048
$ws.Range("Q14").Value = 1862000
$ws.Range("R14").Value = 1841000
$ws.Range("S14").Value = 1819000
$ws.Range("T14").Value = 1865000
$ws.Range("U14").Value = 1875000

# --- Row 29: 5yr, 5yr Forward ---
$ws.Range("N29").Value = 46064
$ws.Range("Q29").Value = 2.15
$ws.Range("R29").Value = 2.17
$ws.Range("S29").Value = 2.2
$ws.Range("T29").Value = 2.18
$ws.Range("U29").Value = 2.16

# --- Row 30: 10yr TIPS ---
$ws.Range("N30").Value = 46064
$ws.Range("Q30").Value = 2.32
$ws.Range("R30").Value = 2.32
$ws.Range("S30").Value = 2.35
$ws.Range("T30").Value = 2.34
$ws.Range("U30").Value = 2.32

# --- Row 47: FFR ---
$ws.Range("N47").Value = 46063

# --- Row 48: 2y UST ---
$ws.Range("N48").Value = 46063
$ws.Range("Q48").Value = 3.45
$ws.Range("R48").Value = 3.48
$ws.Range("S48").Value = 3.5
$ws.Range("T48").Value = 3.47

# --- Row 49: 5y UST ---
$ws.Range("N49").Value = 46063
$ws.Range("Q49").Value = 3.7
$ws.Range("R49").Value = 3.75
$ws.Range("S49").Value = 3.76
$ws.Range("T49").Value = 3.74

# --- Row 50: 10y UST ---
$ws.Range("N50").Value = 46063
$ws.Range("Q50").Value = 4.16
$ws.Range("R50").Value = 4.22
$ws.Range("S50").Value = 4.22
$ws.Range("T50").Value = 4.21
$ws.Range("U50").Value = 4.29

# --- Row 51: 30y Mtg. -> drop the highlight (value unchanged) ---
Clear-DateHighlight("N51")

# --- Row 52: BAA ---
$ws.Range("N52").Value = 46063
$ws.Range("Q52").Value = 5.82
$ws.Range("R52").Value = 5.86
$ws.Range("S52").Value = 5.87
$ws.Range("T52").Value = 5.88
$ws.Range("U52").Value = 5.93
